$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, bordered, centered) onto the new I1:J1 headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I (I0) and J (IF) columns, rows 2..29
$values = @(
    @(6, 6),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(6, 7),
    @(4, 4),
    @(9, 9),
    @(8, 9),
    @(6, 7),
    @(10, 10),
    @(8, 8),
    @(6, 6),
    @(9, 9),
    @(10, 10),
    @(7, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(6, 6),
    @(9, 9),
    @(9, 9),
    @(7, 7),
    @(6, 6),
    @(5, 5),
    @(6, 6),
    @(5, 5),
    @(5, 5)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $pair = $values[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
